$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# OLE color value for RGB(0x16,0x16,0x16), used by the new small bold
# "Noto Sans" header/note font introduced by this edit.
$newFontColor = 1447446
# xlHAlignGeneral / xlHAlignLeft constants (the sheet's column default is
# left-aligned, so "general" must be set explicitly where needed).
$xlHAlignGeneral = 1
$xlHAlignLeft = -4131

# --- Row 5 ---
# G5: the "Sun Compass Information" banner moves up from row 6 to row 5,
# keeping its original bold / Arial / size-10 / general-aligned styling.
$g5 = $ws.Range("G5")
$g5.Value() = "---------------–------------------------------------------------------------------------------Sun Compass Information---------------------------------------------------------------------------------------------"
$g5.Font.Bold() = $true
$g5.HorizontalAlignment() = $xlHAlignGeneral

# R5: new "Calculated Fields" banner, styled with the new bold 9pt Noto Sans
# font (general alignment).
$r5 = $ws.Range("R5")
$r5.Value() = "---------------------------------------------------Calculated Fields----------------------------------------------------"
$r5.Font.Bold() = $true
$r5.Font.Size() = 9
$r5.Font.Name() = "Noto Sans"
$r5.Font.Color() = $newFontColor
$r5.HorizontalAlignment() = $xlHAlignGeneral

# --- Row 6 ---
# E6, G6 and Q6 reuse the exact same new-font / general-alignment
# formatting as R5, so copy R5's format over before writing their text.
$e6 = $ws.Range("E6")
$r5.Copy($e6)
$e6.Value() = "only used if no sun data"

$g6 = $ws.Range("G6")
$r5.Copy($g6)
$g6.Value() = "all sun compass info is optional"

$q6 = $ws.Range("Q6")
$r5.Copy($q6)
$q6.Value() = "Optional Field"

# R6: new note, same new font but left-aligned (unlike the others in this row).
$r6 = $ws.Range("R6")
$r6.Value() = "default core strike"
$r6.Font.Bold() = $true
$r6.Font.Size() = 9
$r6.Font.Name() = "Noto Sans"
$r6.Font.Color() = $newFontColor
$r6.HorizontalAlignment() = $xlHAlignLeft

# --- Sheet view changes ---
# Final selection moves from P6 to Q7 (the scroll position also shifts from
# H1 to G1, but that is not independently observable through this API).
$ws.Range("Q7").Select()
